$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9093090
$ws.Range("J17").Value = 9093090
$ws.Range("L17").Value = 27279270
$ws.Range("N17").Value = -27279606

$ws.Range("H62").Value = 6098.5
$ws.Range("I62").Value = 6289.364
$ws.Range("K62").Value = 6289.364
$ws.Range("M62").Value = -5665.364

$ws.Range("H65").Value = 6098.5
$ws.Range("I65").Value = 6289.364
$ws.Range("K65").Value = 31446.82
$ws.Range("M65").Value = -28326.82

$ws.Range("H86").Value = 2122.4375
$ws.Range("I86").Value = 2823
$ws.Range("J86").Value = 1702.1
$ws.Range("K86").Value = 2823
$ws.Range("L86").Value = 1702.1
$ws.Range("M86").Value = -1700
$ws.Range("N86").Value = -3948.1

$ws.Range("H89").Value = 2122.4375
$ws.Range("I89").Value = 2823
$ws.Range("J89").Value = 1702.1
$ws.Range("K89").Value = 14115
$ws.Range("L89").Value = 8510.5
$ws.Range("M89").Value = -8499
$ws.Range("N89").Value = -19742.5

$ws.Range("H97").Value = 1165.4
$ws.Range("J97").Value = 1165.4
$ws.Range("L97").Value = 3496.2
$ws.Range("N97").Value = -4488.200000000001

$ws.Range("H98").Value = 774.6818
$ws.Range("I98").Value = 775.1177
$ws.Range("K98").Value = 775.1177
$ws.Range("M98").Value = 722.8823

$ws.Range("H122").Value = 774.6818
$ws.Range("I122").Value = 775.1177
$ws.Range("K122").Value = 2325.3531
$ws.Range("M122").Value = 124.6468999999997

$ws.Range("H127").Value = 2999.6667
$ws.Range("I127").Value = 2999
$ws.Range("K127").Value = 8997
$ws.Range("M127").Value = -4037

$ws.Range("H129").Value = 1394.6666
$ws.Range("I129").Value = 926.3333
$ws.Range("K129").Value = 2778.9999
$ws.Range("M129").Value = 2221.0001

$ws.Range("H135").Value = 721.1667
$ws.Range("I135").Value = 747.41174
$ws.Range("K135").Value = 6726.70566
$ws.Range("M135").Value = -4191.70566

$ws.Range("H137").Value = 1289.2
$ws.Range("I137").Value = 1099.7222
$ws.Range("K137").Value = 3299.1666
$ws.Range("M137").Value = -749.1665999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5952.8
$ws.Range("I32").Value = 5115.7856
$ws.Range("K32").Value = 5115.7856
$ws.Range("M32").Value = -4828.7856

$ws.Range("H74").Value = 6211.5293
$ws.Range("I74").Value = 3892.077
$ws.Range("K74").Value = 3892.077
$ws.Range("M74").Value = -3018.077

$ws.Range("H77").Value = 6211.5293
$ws.Range("I77").Value = 3892.077
$ws.Range("K77").Value = 19460.385
$ws.Range("M77").Value = -15092.385

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 695.95
$ws.Range("J94").Value = 1258.875
$ws.Range("L94").Value = 1258.875
$ws.Range("N94").Value = -2160.875

$ws.Range("H105").Value = 3418.8
$ws.Range("I105").Value = 3418.8
$ws.Range("K105").Value = 3418.8
$ws.Range("M105").Value = -1671.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6125.96
$ws.Range("I31").Value = 5405.625
$ws.Range("J31").Value = 6464.9414
$ws.Range("K31").Value = 5405.625
$ws.Range("L31").Value = 6464.9414
$ws.Range("M31").Value = -5110.625
$ws.Range("N31").Value = -7054.9414

$ws.Range("H34").Value = 6125.96
$ws.Range("I34").Value = 5405.625
$ws.Range("J34").Value = 6464.9414
$ws.Range("K34").Value = 5405.625
$ws.Range("L34").Value = 6464.9414
$ws.Range("M34").Value = -5203.625
$ws.Range("N34").Value = -6868.9414

$ws.Range("H82").Value = 52749.25
$ws.Range("J82").Value = 52749.25
$ws.Range("L82").Value = 52749.25
$ws.Range("N82").Value = -53471.25

$ws.Range("H85").Value = 52749.25
$ws.Range("J85").Value = 52749.25
$ws.Range("L85").Value = 52749.25
$ws.Range("N85").Value = -55245.25

$ws.Range("H132").Value = 2450.4285
$ws.Range("I132").Value = 1897.8823
$ws.Range("K132").Value = 5693.6469
$ws.Range("M132").Value = -3163.6469

$ws.Range("H141").Value = 316996.56
$ws.Range("J141").Value = 316996.56
$ws.Range("L141").Value = 316996.56
$ws.Range("N141").Value = -327356.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22948158
$ws.Range("I4").Value = 27391064
$ws.Range("K4").Value = 82173192
$ws.Range("M4").Value = -82173080

$ws.Range("H39").Value = 394.25
$ws.Range("I39").Value = 384
$ws.Range("K39").Value = 1152
$ws.Range("M39").Value = -858

$ws.Range("H55").Value = 871.7
$ws.Range("J55").Value = 912.8333
$ws.Range("L55").Value = 2738.4999
$ws.Range("N55").Value = -3092.4999

$ws.Range("H131").Value = 15627469
$ws.Range("J131").Value = 2808.3462
$ws.Range("L131").Value = 8425.0386
$ws.Range("N131").Value = -18505.0386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4488.0835
$ws.Range("I102").Value = 2826.6316
$ws.Range("K102").Value = 2826.6316
$ws.Range("M102").Value = -1204.6316

$ws.Range("H113").Value = 335660.9
$ws.Range("I113").Value = 572928.3
$ws.Range("K113").Value = 572928.3
$ws.Range("M113").Value = -570758.3

$ws.Range("H122").Value = 2151.9092
$ws.Range("I122").Value = 2117.2
$ws.Range("K122").Value = 6351.599999999999
$ws.Range("M122").Value = -3901.599999999999

$ws.Range("H132").Value = 4042.0952
$ws.Range("J132").Value = 4200.25
$ws.Range("L132").Value = 12600.75
$ws.Range("N132").Value = -17660.75

$ws.Range("H134").Value = 73743.914
$ws.Range("J134").Value = 73743.914
$ws.Range("L134").Value = 221231.742
$ws.Range("N134").Value = -226301.742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5947
$ws.Range("J7").Value = 5663.75
$ws.Range("L7").Value = 5663.75
$ws.Range("N7").Value = -5887.75

$ws.Range("H126").Value = 5947
$ws.Range("J126").Value = 5663.75
$ws.Range("L126").Value = 16991.25
$ws.Range("N126").Value = -21931.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 15000
$ws.Range("I52").Value = 15000
$ws.Range("K52").Value = 15000
$ws.Range("M52").Value = -14774

$ws.Range("H62").Value = 3199.7144
$ws.Range("I62").Value = 3066
$ws.Range("J62").Value = 3300
$ws.Range("K62").Value = 3066
$ws.Range("L62").Value = 3300
$ws.Range("M62").Value = -2442
$ws.Range("N62").Value = -4548

$ws.Range("H65").Value = 3199.7144
$ws.Range("I65").Value = 3066
$ws.Range("J65").Value = 3300
$ws.Range("K65").Value = 15330
$ws.Range("L65").Value = 16500
$ws.Range("M65").Value = -12210
$ws.Range("N65").Value = -22740

$ws.Range("H100").Value = 520.8095
$ws.Range("I100").Value = 451.88235
$ws.Range("J100").Value = 813.75
$ws.Range("K100").Value = 903.7646999999999
$ws.Range("L100").Value = 1627.5
$ws.Range("M100").Value = -362.7646999999999
$ws.Range("N100").Value = -2709.5

$ws.Range("H107").Value = 1226.8
$ws.Range("I107").Value = 1059.8572
$ws.Range("J107").Value = 1616.3334
$ws.Range("K107").Value = 3179.5716
$ws.Range("L107").Value = 4849.0002
$ws.Range("M107").Value = -1259.5716
$ws.Range("N107").Value = -8689.0002

$ws.Range("H122").Value = 3547.4546
$ws.Range("I122").Value = 2324.5715
$ws.Range("K122").Value = 6973.7145
$ws.Range("M122").Value = -4523.7145

$ws.Range("H132").Value = 5074.778
$ws.Range("I132").Value = 4250.8335
$ws.Range("K132").Value = 12752.5005
$ws.Range("M132").Value = -10222.5005
